# Update the cryptos price list with the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> (new Price, new Volume(1h)) values.
# Only rows whose values changed are listed (rows 4 and 37 are unchanged).
$updates = @(
    @{ Row = 2;  D = "64.136.71";  E = "  +0.56%  " },
    @{ Row = 3;  D = "3.148.50";   E = "  +0.45%  " },
    @{ Row = 5;  D = "591.89";     E = "  +0.50%  " },
    @{ Row = 6;  D = "146.21";     E = "  -0.78%  " },
    @{ Row = 7;  D = $null;        E = "  +0.03%  " },
    @{ Row = 8;  D = "3.139.09";   E = "  +0.34%  " },
    @{ Row = 9;  D = $null;        E = "  -0.76%  " },
    @{ Row = 10; D = $null;        E = "  +0.50%  " },
    @{ Row = 11; D = "5.83";       E = "  +1.81%  " },
    @{ Row = 12; D = $null;        E = "  -2.33%  " },
    @{ Row = 13; D = $null;        E = "  -2.61%  " },
    @{ Row = 14; D = $null;        E = "  -0.45%  " },
    @{ Row = 15; D = "3.665.16";   E = "  +0.33%  " },
    @{ Row = 16; D = $null;        E = "  -1.25%  " },
    @{ Row = 17; D = "7.26";       E = "  +1.03%  " },
    @{ Row = 18; D = "63.942.99";  E = "  +0.43%  " },
    @{ Row = 19; D = "3.143.90";   E = "  +0.33%  " },
    @{ Row = 20; D = "466.64";     E = "  +0.15%  " },
    @{ Row = 21; D = "14.37";      E = "  -0.21%  " },
    @{ Row = 22; D = "0.732";      E = "  -0.30%  " },
    @{ Row = 23; D = $null;        E = "  -0.87%  " },
    @{ Row = 24; D = "13.02";      E = "  -2.01%  " },
    @{ Row = 25; D = "81.26";      E = "  -1.25%  " },
    @{ Row = 26; D = $null;        E = "  +5.82%  " },
    @{ Row = 27; D = $null;        E = "  +0.07%  " },
    @{ Row = 28; D = "9.64";       E = "  +7.48%  " },
    @{ Row = 29; D = $null;        E = "  +1.14%  " },
    @{ Row = 30; D = $null;        E = "  -0.04%  " },
    @{ Row = 31; D = $null;        E = "  +6.69%  " },
    @{ Row = 32; D = $null;        E = "  +0.11%  " },
    @{ Row = 33; D = "27.52";      E = "  +1.12%  " },
    @{ Row = 34; D = $null;        E = "  +0.95%  " },
    @{ Row = 35; D = $null;        E = "  -5.80%  " },
    @{ Row = 36; D = "1.06";       E = $null },
    @{ Row = 38; D = "2.30";       E = "  -2.86%  " },
    @{ Row = 39; D = "3.23";       E = $null },
    @{ Row = 40; D = "459.10";     E = "  +0.83%  " },
    @{ Row = 41; D = "51.31";      E = "  +0.49%  " },
    @{ Row = 42; D = "9.24";       E = "  +5.44%  " },
    @{ Row = 43; D = "0.293";      E = "  +4.64%  " },
    @{ Row = 44; D = $null;        E = "  -0.19%  " },
    @{ Row = 45; D = "2.925.13";   E = "  +0.67%  " },
    @{ Row = 46; D = "39.55";      E = "  +10.07%  " },
    @{ Row = 47; D = $null;        E = "  -2.52%  " },
    @{ Row = 48; D = "131.77";     E = "  +3.03%  " },
    @{ Row = 49; D = $null;        E = "  -0.06%  " },
    @{ Row = 50; D = "2.26";       E = "  +2.80%  " },
    @{ Row = 51; D = $null;        E = "  -0.97%  " }
)

# Rows whose Price text (column D) parses as a plain decimal number (e.g.
# "591.89") and therefore needs to be forced into the cell as literal text
# -- otherwise Excel auto-converts it to a numeric value. Rows using the
# dotted "thousands" style (e.g. "64.136.71") are never ambiguous and don't
# need this treatment.
$numericLooking = @(5, 6, 11, 17, 20, 21, 22, 24, 25, 28, 33, 36, 38, 39, 40, 41, 42, 43, 46, 48, 50)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($numericLooking -contains $r) {
            # Leading apostrophe forces text entry; ClearFormats() then
            # drops the "quote prefix" number format Excel applies so the
            # cell's style matches the rest of the (untouched) column.
            $cell.Value = "'" + $u.D
            $cell.ClearFormats()
        } else {
            $cell.Value = $u.D
        }
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
